$d = $word.ActiveDocument

# Prepend "PD " to the "Dr. med. Thiên-Trí Lâm" line in the main body
# (the title is being updated to reflect the "Privatdozent" academic title).
$d.Content.Find.Execute("Dr. med. Thiên-Trí Lâm", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PD Dr. med. Thiên-Trí Lâm", 2)
